# Add the "Github: https://github.com/beken25/yapayzekafinal" line to the
# empty trailing paragraph of the document (the paragraph immediately
# before the final sectPr).
#
# The target markup is:
#   <w:p>
#     <w:proofErr w:type="spellStart"/>
#     <w:r><w:t>Github</w:t></w:r>
#     <w:proofErr w:type="spellEnd"/>
#     <w:r><w:t xml:space="preserve">: </w:t></w:r>
#     <w:r><w:t>https://github.com/beken25/yapayzekafinal</w:t></w:r>
#   </w:p>
#
# Plain text assignment (Range.Text = "...") would merge everything into a
# single run and would not reproduce the proofErr spell-check markers, so
# the paragraph's content is inserted as literal WordprocessingML via
# Range.InsertXML - this lets us create the exact run/proofErr layout the
# diff expects. Inserting at a collapsed range positioned at the start of
# the existing (empty) last paragraph adds the content *into* that
# paragraph instead of creating a brand new one, so the paragraph count
# does not change (matches the diff, which turns "<w:p/>" into a "<w:p>"
# that merely gained children).

$d = $word.ActiveDocument

# The last paragraph in the document body is the empty one that sits right
# before the sectPr - that's the one the diff turns from "<w:p/>" into the
# paragraph holding the Github link.
$target = $d.Paragraphs.Last

# Collapsed range at the very start of that (empty) paragraph.
$insertionPoint = $d.Range($target.Range.Start, $target.Range.Start)

$wordNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$bodyXml = '<w:p xmlns:w="' + $wordNs + '">' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Github</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">: </w:t></w:r>' +
    '<w:r><w:t>https://github.com/beken25/yapayzekafinal</w:t></w:r>' +
    '</w:p>'

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' + $bodyXml + '</pkg:xmlData>' +
    '</pkg:part>' +
    '</pkg:package>'

$insertionPoint.InsertXML($packageXml)
